$wb = $excel.ActiveWorkbook

# Both "展览" (sheet1) and "全部类型" (sheet4) carry the same event rows,
# so the "想去人数" (column F) / "最低票价" (column G) refresh applies to both.
foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Cells.Item(2, 6).Value = 197
    $ws.Cells.Item(3, 6).Value = 435
    $ws.Cells.Item(4, 6).Value = 12705
    $ws.Cells.Item(5, 6).Value = 1293
    $ws.Cells.Item(5, 7).Value = 1
    $ws.Cells.Item(6, 6).Value = 162
    $ws.Cells.Item(8, 6).Value = 93
    $ws.Cells.Item(9, 6).Value = 162
    $ws.Cells.Item(10, 6).Value = 212
    $ws.Cells.Item(11, 6).Value = 459
    $ws.Cells.Item(12, 6).Value = 61
    $ws.Cells.Item(16, 6).Value = 387
    $ws.Cells.Item(17, 6).Value = 5450
    $ws.Cells.Item(19, 6).Value = 27
    $ws.Cells.Item(20, 6).Value = 949
    $ws.Cells.Item(21, 6).Value = 25
    $ws.Cells.Item(22, 6).Value = 128
    $ws.Cells.Item(23, 6).Value = 90
}
